$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "44.137.96"
$ws.Range("E2").Value = "  +3.37%  "
$ws.Range("D3").Value = "2.250.74"
$ws.Range("E3").Value = "  +2.17%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "'258.22"
$ws.Range("E5").Value = "  +2.64%  "
$ws.Range("D6").Value = "'80.93"
$ws.Range("E6").Value = "  +8.22%  "
$ws.Range("E7").Value = "  +2.57%  "
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("D9").Value = "'0.606"
$ws.Range("E9").Value = "  +2.88%  "
$ws.Range("D10").Value = "'43.61"
$ws.Range("E10").Value = "  +7.92%  "
$ws.Range("D11").Value = "'0.0935"
$ws.Range("E11").Value = "  +1.31%  "
$ws.Range("D12").Value = "'7.12"
$ws.Range("E12").Value = "  +4.18%  "
$ws.Range("E13").Value = "  +2.17%  "
$ws.Range("D14").Value = "2.585.66"
$ws.Range("E14").Value = "  +2.44%  "
$ws.Range("D15").Value = "'14.83"
$ws.Range("E15").Value = "  +3.20%  "
$ws.Range("D16").Value = "2.252.98"
$ws.Range("E16").Value = "  +2.60%  "
$ws.Range("D17").Value = "'0.798"
$ws.Range("E17").Value = "  +1.92%  "
$ws.Range("D18").Value = "44.068.18"
$ws.Range("E18").Value = "  +3.52%  "
$ws.Range("E19").Value = "  +2.40%  "
$ws.Range("D20").Value = "'71.74"
$ws.Range("E20").Value = "  +0.68%  "
$ws.Range("E21").Value = "  +2.69%  "
$ws.Range("D22").Value = "'2.36"
$ws.Range("E22").Value = "  +9.52%  "
$ws.Range("D23").Value = "'235.41"
$ws.Range("E23").Value = "  +2.76%  "
$ws.Range("D24").Value = "'9.50"
$ws.Range("E24").Value = "  +0.75%  "
$ws.Range("E25").Value = "  +0.15%  "
$ws.Range("D26").Value = "'10.92"
$ws.Range("E26").Value = "  +1.78%  "
$ws.Range("D27").Value = "'41.17"
$ws.Range("E27").Value = "  +8.27%  "
$ws.Range("D28").Value = "'3.36"
$ws.Range("E28").Value = "  -0.46%  "
$ws.Range("D29").Value = "'2.25"
$ws.Range("E29").Value = "  +1.66%  "
$ws.Range("D30").Value = "'2.21"
$ws.Range("E30").Value = "  -0.67%  "
$ws.Range("D31").Value = "'173.34"
$ws.Range("E31").Value = "  +2.32%  "
$ws.Range("E32").Value = "  +2.77%  "
$ws.Range("D33").Value = "'0.0878"
$ws.Range("E33").Value = "  +9.54%  "
$ws.Range("E34").Value = "  +3.25%  "
$ws.Range("D35").Value = "'0.115"
$ws.Range("E35").Value = "  +6.41%  "
$ws.Range("D36").Value = "'0.123"
$ws.Range("E36").Value = "  +1.81%  "
$ws.Range("D37").Value = "'0.0369"
$ws.Range("E37").Value = "  +12.54%  "
$ws.Range("D38").Value = "'4.53"
$ws.Range("E38").Value = "  +3.75%  "
$ws.Range("D39").Value = "'13.09"
$ws.Range("E39").Value = "  +6.53%  "
$ws.Range("D40").Value = "'2.94"
$ws.Range("E40").Value = "  +21.59%  "
$ws.Range("E41").Value = "  +3.62%  "
$ws.Range("D42").Value = "'63.43"
$ws.Range("E42").Value = "  +7.16%  "
$ws.Range("D43").Value = "'5.55"
$ws.Range("E43").Value = "  +5.27%  "
$ws.Range("E44").Value = "  +2.89%  "
$ws.Range("D45").Value = "'104.31"
$ws.Range("E45").Value = "  +0.85%  "
$ws.Range("D46").Value = "'8.58"
$ws.Range("E46").Value = "  +1.27%  "
$ws.Range("E47").Value = "  +2.00%  "
$ws.Range("D48").Value = "'0.465"
$ws.Range("E48").Value = "  -3.48%  "
$ws.Range("D49").Value = "'1.12"
$ws.Range("E49").Value = "  +1.85%  "
$ws.Range("B50").Value = "Stacks"
$ws.Range("C50").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D50").Value = "'1.54"
$ws.Range("E50").Value = "  +25.75%  "
$ws.Range("B51").Value = "TrustWalletToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D51").Value = "'1.15"
$ws.Range("E51").Value = "  +2.13%  "
